# "new .ttl from Google sheet has been generated"
#
# The regenerated export:
#   - upgrades two existing dct:creator ORCID URLs (rows 10-11) from
#     http:// to https://
#   - gains one additional dct:creator row (new row 12, ORCID
#     0000-0003-0732-4617) that wasn't present before, which pushes every
#     row from the old row 12 onward down by one (old row 85 -> new row 86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix scheme on the two pre-existing creator rows.
$ws.Range("B10").Value = "https://orcid.org/0000-0003-4124-9040"
$ws.Range("B11").Value = "https://orcid.org/0000-0002-9381-9693"

# Insert the new creator row; this shifts rows 12:85 down to 13:86 and
# bumps the sheet dimension from A1:S85 to A1:S86 automatically.
$ws.Rows(12).Insert()

$ws.Range("A12").Value = "dct:creator"
$ws.Range("B12").Value = "https://orcid.org/0000-0003-0732-4617"
